$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -12.8112
$ws.Range("D4").Value = -8.196400000000004
$ws.Range("C6").Value = -11.8075
$ws.Range("C7").Value = -12.7054
$ws.Range("D9").Value = -7.427299999999991
$ws.Range("D12").Value = -5.777999999999999
$ws.Range("C16").Value = -14.58959999999999
$ws.Range("D17").Value = -8.615799999999991
$ws.Range("D18").Value = -8.864999999999995
$ws.Range("D19").Value = -8.418599999999994
$ws.Range("C20").Value = -11.829
$ws.Range("D20").Value = -8.483799999999997
$ws.Range("D26").Value = -7.815099999999997
$ws.Range("C28").Value = -12.9487
$ws.Range("C29").Value = -11.75870000000001
$ws.Range("D31").Value = -7.403699999999992
$ws.Range("C32").Value = -12.5565
$ws.Range("D39").Value = -8.320899999999988
$ws.Range("C40").Value = -11.6697
$ws.Range("D40").Value = -8.132299999999994
$ws.Range("D41").Value = -7.567299999999999
$ws.Range("D42").Value = -7.840399999999997
$ws.Range("D43").Value = -7.186100000000005
$ws.Range("C46").Value = -14.72809999999999
$ws.Range("D47").Value = -7.563299999999999
$ws.Range("D48").Value = -7.205499999999997
$ws.Range("C51").Value = -11.6097
$ws.Range("C52").Value = -11.5241
$ws.Range("C57").Value = -14.4439
$ws.Range("C59").Value = -12.96760000000001
$ws.Range("C62").Value = -14.49969999999999
$ws.Range("D63").Value = -6.610699999999993
$ws.Range("D64").Value = -6.931299999999994
$ws.Range("C66").Value = -11.6609
$ws.Range("C73").Value = -10.99660000000001
$ws.Range("C74").Value = -12.03580000000001
$ws.Range("D76").Value = -7.4958
$ws.Range("D81").Value = -7.620199999999999
$ws.Range("D89").Value = -8.296399999999998
$ws.Range("C92").Value = -10.4024
$ws.Range("D94").Value = -5.883499999999998
$ws.Range("C100").Value = -11.6673
